# Apply updated crypto price/volume figures (cryptos list refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a plain number (e.g. "250.29") need a
# leading apostrophe so Excel stores them as text, matching the source data
# (prices such as "37.150.81" use dots as thousands separators, not as a
# decimal point, so they must stay text rather than become real numbers).

$ws.Range("D2").Value = "37.150.81"
$ws.Range("E2").Value = "  +0.17%  "
$ws.Range("D3").Value = "2.056.02"
$ws.Range("E3").Value = "  -0.14%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "'250.29"
$ws.Range("E5").Value = "  +0.56%  "
$ws.Range("E6").Value = "  +1.79%  "
$ws.Range("D7").Value = "'61.71"
$ws.Range("E7").Value = "  +11.50%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").Value = "'0.387"
$ws.Range("E10").Value = "  -0.76%  "
$ws.Range("E11").Value = "  +2.03%  "
$ws.Range("D12").Value = "'16.35"
$ws.Range("E12").Value = "  +8.14%  "
$ws.Range("D13").Value = "2.356.18"
$ws.Range("E13").Value = "  -0.06%  "
$ws.Range("E14").Value = "  +1.27%  "
$ws.Range("D15").Value = "'5.73"
$ws.Range("E15").Value = "  +9.24%  "
$ws.Range("D16").Value = "2.054.88"
$ws.Range("E16").Value = "  -0.15%  "
$ws.Range("D17").Value = "'18.02"
$ws.Range("E17").Value = "  +26.92%  "
$ws.Range("D18").Value = "37.150.07"
$ws.Range("E18").Value = "  +0.25%  "
$ws.Range("D19").Value = "'75.28"
$ws.Range("E19").Value = "  +3.87%  "
$ws.Range("E20").Value = "  -4.32%  "
$ws.Range("E21").Value = "  +1.29%  "
$ws.Range("D22").Value = "'239.41"
$ws.Range("E22").Value = "  +0.90%  "
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("E24").Value = "  -1.06%  "
$ws.Range("D25").Value = "'2.21"
$ws.Range("E25").Value = "  +12.15%  "
$ws.Range("D26").Value = "'169.49"
$ws.Range("E26").Value = "  -0.76%  "
$ws.Range("D27").Value = "'9.44"
$ws.Range("E27").Value = "  +4.27%  "
$ws.Range("D28").Value = "'20.04"
$ws.Range("E28").Value = "  -0.65%  "
$ws.Range("E29").Value = "  +1.82%  "
$ws.Range("E30").Value = "  +9.87%  "
$ws.Range("E31").Value = "  +5.75%  "
$ws.Range("D32").Value = "'0.0622"
$ws.Range("E32").Value = "  -0.43%  "
$ws.Range("E33").Value = "  +4.30%  "
$ws.Range("E34").Value = "  +4.90%  "
$ws.Range("E35").Value = "  -0.26%  "
$ws.Range("E36").Value = "  -1.13%  "
$ws.Range("E37").Value = "  -1.27%  "
$ws.Range("E38").Value = "  +2.24%  "
$ws.Range("E39").Value = "  +0.74%  "
$ws.Range("D40").Value = "'5.32"
$ws.Range("E40").Value = "  +32.27%  "
$ws.Range("D41").Value = "'3.18"
$ws.Range("E41").Value = "  +14.51%  "
$ws.Range("D42").Value = "'18.19"
$ws.Range("E42").Value = "  +0.48%  "
$ws.Range("E43").Value = "  +0.43%  "
$ws.Range("D44").Value = "'98.16"
$ws.Range("E44").Value = "  +1.76%  "
$ws.Range("E45").Value = "  -0.09%  "
$ws.Range("E46").Value = "  +2.57%  "
$ws.Range("D47").Value = "1.297.46"
$ws.Range("E47").Value = "  +0.04%  "
$ws.Range("E48").Value = "  -1.43%  "
$ws.Range("D49").Value = "'6.86"
$ws.Range("E49").Value = "  +0.92%  "
$ws.Range("D50").Value = "2.242.50"
$ws.Range("E50").Value = "  -0.44%  "
$ws.Range("D51").Value = "'3.56"
$ws.Range("E51").Value = "  -15.74%  "
